# Insert a new weekly price record for "Pepino dulce" (Vega Modelo de Temuco)
# as row 66, pushing the previously-existing rows 66..159 down to 67..160.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 66..159 down by inserting a brand-new row 66.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record's values.
$ws.Range("A66").Value = 10
$ws.Range("B66").Value = "Vega Modelo de Temuco"
$ws.Range("C66").Value = "La Araucanía"
$ws.Range("D66").Value = 44482
$ws.Range("E66").Value = 9
$ws.Range("F66").Value = 100112043
$ws.Range("G66").Value = "Pepino dulce"
$ws.Range("H66").Value = "Cultivar IV Región"
$ws.Range("I66").Value = "Segunda"
$ws.Range("J66").Value = 50
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = 20000
$ws.Range("N66").Value = "`$/bandeja 18 kilos"
$ws.Range("O66").Value = "Provincia de Limarí"
$ws.Range("P66").Value = 1111
$ws.Range("Q66").Value = 18
$ws.Range("R66").Value = "Hortaliza"
